$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# Column M = "TieneVencimiento" (boolean): set these rows to TRUE
$mRows = @(2,3,4,5,7,8,13,17,18,19,20,21)
foreach ($r in $mRows) {
    $ws.Range("M$r").Value = $true
}

# Column O = "ImagenExactaDelArticulo" (boolean): set these rows to TRUE
$oRows = @(22,25)
foreach ($r in $oRows) {
    $ws.Range("O$r").Value = $true
}
